$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new DSLOOKUP demo formula in A6 (row 5 is intentionally left empty)
$ws.Range("A6").Formula = '=DSLOOKUP("pers", "ADDRESS", "Abc St.", "CITY", "Lublin", "LASTNAME")'

# Move the active selection to D9 (was G12)
$ws.Range("D9").Select()

# Reflect the workbook window's new vertical screen position (yWindow 5505 -> 6405)
$excel.Windows.Item(1).Top = 6405
